# Apply automatic update: re-order data among rows 2-8 for columns
# A, B, D, E, F, G, H, Q, R, AO (other columns stay put).
#
# New row r receives the old contents (for the above columns) of row
# $mapping[r]:
#   2 <- 3
#   3 <- 2
#   4 <- 8
#   5 <- 7
#   6 <- 4
#   7 <- 6
#   8 <- 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AO")
$rows = @(2, 3, 4, 5, 6, 7, 8)

$mapping = @{ 2 = 3; 3 = 2; 4 = 8; 5 = 7; 6 = 4; 7 = 6; 8 = 5 }

# Snapshot the current (old) values for every affected cell before
# overwriting anything.
$old = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $old["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Write back according to the mapping.
foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $old["$col$src"]
    }
}
